$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet is protected; unprotect (no known password required by this host)
# so the cell values below can be written, then restore protection afterward.
$ws.Unprotect("")

# Update the "as of" date in the confidential disclaimer banner (A59).
$disclaimer = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-05-10 for illustrative purposes only and are subject to change."
$ws.Range("A59").Value = $disclaimer

# Refresh Weight (D) and Percent Change (E) figures for each holding row.
$ws.Range("D2").Value = 0.01300437941203115
$ws.Range("E2").Value = -0.01106259097525464
$ws.Range("D3").Value = 0.01061897222406256
$ws.Range("E3").Value = -0.0106343967730107
$ws.Range("D4").Value = 0.01041258955890843
$ws.Range("E4").Value = -0.002493143854400359
$ws.Range("D5").Value = 0.0114142161204007
$ws.Range("E5").Value = -0.02467685076380721
$ws.Range("D6").Value = 0.01087467904818435
$ws.Range("E6").Value = -0.006922893291955057
$ws.Range("D7").Value = 0.01354867583291982
$ws.Range("E7").Value = 0.007704160246533309
$ws.Range("D8").Value = 0.010918594856388
$ws.Range("E8").Value = -0.01462225832656372
$ws.Range("D9").Value = 0.01114552925444533
$ws.Range("E9").Value = -0.01824534161490687
$ws.Range("D10").Value = 0.01041648357145851
$ws.Range("E10").Value = -0.01781931464174469
$ws.Range("D11").Value = 0.01106224065267977
$ws.Range("E11").Value = -0.02156057494866548
$ws.Range("D12").Value = 0.442606446483693
$ws.Range("E12").Value = -0.003783102143757766
$ws.Range("D13").Value = 0.01162881947871611
$ws.Range("E13").Value = -0.005692599620493288
$ws.Range("D14").Value = 0.0106548836731355
$ws.Range("E14").Value = -0.0008527572484366974
$ws.Range("D15").Value = 0.01008030048796845
$ws.Range("E15").Value = -0.005193578847969671
$ws.Range("D16").Value = 0.009907882265612216
$ws.Range("E16").Value = -0.009083167754754484
$ws.Range("D17").Value = 0.009465154672071411
$ws.Range("E17").Value = 0.01508485229415446
$ws.Range("D18").Value = 0.008030751882444087
$ws.Range("E18").Value = -0.02168525402726162
$ws.Range("D19").Value = 0.009015396222537407
$ws.Range("E19").Value = -0.003167480533192379
$ws.Range("D20").Value = 0.01094217526571903
$ws.Range("E20").Value = -0.05634638196915775
$ws.Range("D21").Value = 0.01209285597426707
$ws.Range("E21").Value = -0.0001610046691354894
$ws.Range("D22").Value = 0.01156186409625782
$ws.Range("E22").Value = -0.001665278934221193
$ws.Range("D23").Value = 0.01149544954887594
$ws.Range("E23").Value = -0.02582921665490479
$ws.Range("D24").Value = 0.01261335565179415
$ws.Range("E24").Value = 0.01320641454420723
$ws.Range("D25").Value = 0.01273428637487713
$ws.Range("E25").Value = -0.008969828757814602
$ws.Range("D26").Value = 0.01190010235303821
$ws.Range("E26").Value = -0.01248909249563701
$ws.Range("D27").Value = 0.01258350155557689
$ws.Range("E27").Value = -0.02491103202846978
$ws.Range("D28").Value = 0.01336425107186751
$ws.Range("E28").Value = 0.06837606837606836
$ws.Range("D29").Value = 0.01180210303719458
$ws.Range("E29").Value = -0.05966455870222698
$ws.Range("D30").Value = 0.006911331441311956
$ws.Range("E30").Value = 0.00142421159715167
$ws.Range("D31").Value = 0.004961350573352784
$ws.Range("E31").Value = -0.03395650515070581
$ws.Range("D32").Value = 0.009492088258876116
$ws.Range("E32").Value = -0.01318458417849888
$ws.Range("D33").Value = 0.01210821569043682
$ws.Range("E33").Value = -0.1476683937823836
$ws.Range("D34").Value = 0.01036185922874214
$ws.Range("E34").Value = -0.07635054021608645
$ws.Range("D35").Value = 0.00938911326032961
$ws.Range("E35").Value = 0.003410059676044463
$ws.Range("D36").Value = 0.009670563834088021
$ws.Range("E36").Value = -0.02288488210818296
$ws.Range("D37").Value = 0.01030485521168961
$ws.Range("E37").Value = 0.02516059957173455
$ws.Range("D38").Value = 0.01141053844188119
$ws.Range("E38").Value = 0.001042752867570274
$ws.Range("D39").Value = 0.01413613089290519
$ws.Range("E39").Value = -0.01334476003917728
$ws.Range("D40").Value = 0.01139388072152807
$ws.Range("E40").Value = -0.03751803751803739
$ws.Range("D41").Value = 0.0126044859565412
$ws.Range("E41").Value = 0.007963751201427893
$ws.Range("D42").Value = 0.0114617014401086
$ws.Range("E42").Value = -0.0213376367222522
$ws.Range("D43").Value = 0.01151210726922905
$ws.Range("E43").Value = -0.01400933955970651
$ws.Range("D44").Value = 0.01075991384497233
$ws.Range("E44").Value = -0.0268208092485549
$ws.Range("D45").Value = 0.01141973263817998
$ws.Range("E45").Value = -0.01089273028652615
$ws.Range("D46").Value = 0.01109880110384439
$ws.Range("E46").Value = 0.01693825042881647
$ws.Range("D47").Value = 0.01003324783632167
$ws.Range("E47").Value = -0.01494226850803726
$ws.Range("D48").Value = 0.009410530329355038
$ws.Range("E48").Value = -0.02758620689655178
$ws.Range("D49").Value = 0.009851851751697206
$ws.Range("E49").Value = -0.01818181818181808
$ws.Range("D50").Value = 0.009295007957036061
$ws.Range("E50").Value = 0.04650188521156262
$ws.Range("D51").Value = 0.009264396691711837
$ws.Range("E51").Value = -0.02632838678793692
$ws.Range("D52").Value = 0.0102390896663994
$ws.Range("E52").Value = -0.02810057046270853
$ws.Range("D53").Value = 0.008767585590531058
$ws.Range("E53").Value = -0.04656040268456363
$ws.Range("D54").Value = 0.004181195975646196
$ws.Range("E54").Value = 0.0007760962359331103
$ws.Range("D55").Value = 0.004064483766159138
$ws.Range("E55").Value = 0.008090270385352483
$ws.Range("E56").Value = -0.009396955368937365

# Restore sheet protection.
$ws.Protect("")
